# Update EE-2020 abstract: new title/bookmark, author-affiliation line breaks,
# and revised abstract body text.

$d = $word.ActiveDocument

# --- 1. Title heading + its bookmark -------------------------------------
# Replace the whole Heading1 paragraph's XML so the bookmark name and the
# visible title text are updated together, keeping bookmarkStart/bookmarkEnd
# adjacent (zero-width bookmark) exactly as in the original structure.
$titleParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:bookmarkStart w:id="21" w:name="optimally-allocating-resources-for-gathering-evidence-and-managing-biodiversity"/><w:bookmarkEnd w:id="21"/><w:r><w:t xml:space="preserve">Optimally allocating resources for gathering evidence and managing biodiversity</w:t></w:r></w:p>'
$d.Paragraphs(1).Range.InsertXML($titleParaXml)

# --- 2. Author/affiliation paragraph: spaces -> line breaks --------------
# The two plain-space runs that separated the affiliation superscripts from
# the following text become explicit text-wrapping line breaks.
$authorParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">Jeffrey O. Hanson</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t xml:space="preserve">*1</w:t></w:r><w:r><w:t xml:space="preserve">, Joseph R. Bennett</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t xml:space="preserve">2</w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t xml:space="preserve">1</w:t></w:r><w:r><w:t xml:space="preserve">Department of Biology, Carleton University, Ottawa, Canada</w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t xml:space="preserve">*</w:t></w:r><w:r><w:t xml:space="preserve">Presenting author, correspondence should be addressed to jeffrey.hanson@uqconnect.edu.au</w:t></w:r></w:p>'
$d.Paragraphs(2).Range.InsertXML($authorParaXml)

# --- 3. Abstract body text -------------------------------------------------
# Swap the old abstract paragraph for the revised wording via Find/Replace
# (pure text substitution; no structural/formatting change needed here).
$oldBody = 'Conserving biodiversity means working with limited resources and incomplete information. In the face of uncertainty, practitioners can develop plans for expanding protected area systems (prioritisations) using existing data or they can gather evidence -- by surveying sites for species of conservation interest -- to refine them. However, such evidence gathering reduces the funds available for purchasing sites for conservation and so survey schemes need to be strategically designed. Here, we investigated various approaches for generating survey schemes. Using a case-study, we obtained existing survey data for native species, survey costs, and land acquisition costs. Next, we designated new sites and modelled the probability that the study species occupied them. We then used conventional approaches to generate survey schemes by selecting sites with (i) geographically representative locations; (ii) environmentally representative conditions; (iii) uncertain model predictions; (iv) high occupancy probabilities; and (v) low acquisition costs. To compare them, we also generated survey schemes by (vi) directly maximizing return on investment. After generating these survey schemes, we evaluated them using value of information analyses. We found that survey schemes generated by maximizing return on investment were far more effective than conventional approaches. In particular, survey schemes generated by increasing the geographic coverage and environmental diversity of surveyed sites had the poorest performance. Under limited budgets, survey schemes generated using conventional approaches misallocated a large proportion of the available funds so that little remained for purchasing sites to achieve conservation objectives. It was only under relatively large budgets, when the majority of sites could be purchased for conservation, that conventional approaches had near-optimal performance. Our results show that schemes for gathering evidence can be substantially improved by explicitly quantifying their capacity to improve conservation decisions. We recommend using value of information analyses, when feasible, to optimally allocate resources for gathering evidence and conserving biodiversity.'
$newBody = 'Conserving biodiversity means working with limited resources and incomplete information. In the face of uncertainty, practitioners can develop management plans using existing data or they can gather evidence to refine them. Since gathering evidence consumes limited resources, plans for gathering evidence must strategically maximize the amount of information gained to improve management decisions and also ensure that sufficient resources remain for implementing management actions. Here, we investigated different approaches for developing plans to gather evidence. We evaluated these approaches using value of information analyses and decision support tools to simulate management decisions based on different outcomes. We found that directly maximizing return on investment was by far the most effective approach for developing plans to gather evidence. Under limited budgets, alternative approaches produced plans that were highly ineffective. They allocated a large proportion of the available resources towards gathering evidence, so that few resources remained for actually achieving conservation objectives. It was only under relatively large budgets -- when a large amount of resources remained for achieving conservation objectives -- that alternative approaches had near-optimal performance. Our results show that plans for gathering evidence can be substantially improved by explicitly quantifying their capacity to improve conservation decisions. We recommend using value of information analyses, when feasible, to optimally balance the allocation of resources for gathering evidence and conserving biodiversity.'
$d.Content.Find.Execute($oldBody, $true, $false, $false, $false, $false, $true, 1, $false, $newBody, 2)

Write-Output "done"
